$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the orphan "18S ribosomal RNA, Rn18S" row (row 2, which has no
# primer data) - deleting the entire row shifts everything below it up by
# one and drops the now-unused shared string automatically.
$ws.Rows(2).Delete()

# Leave the selection where the editor last clicked before saving.
$ws.Range("D6").Select() | Out-Null
